# Apply the data corrections to the train2Block1Test sheet.
# Each row's audio paths move from "trainaudio/" to "trainingaudio/"
# (and several entries are swapped out for different files), with a
# handful of truePos/falsePos sign flips to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("trainingaudio/18_popata2.wav", "pngimages/18_donut.png", "trainingaudio/06_titoka3.wav", "pngimages/06_tent.png", 0.5, -0.5),
    @("trainingaudio/20_tatito1.wav", "pngimages/20_pizza.png", "trainingaudio/19_papipi1.wav", "pngimages/19_burger.png", -0.5, 0.5),
    @("trainingaudio/21_papika1.wav", "pngimages/21_cheese.png", "trainingaudio/01_kitipi1.wav", "pngimages/01_gift.png", -0.5, 0.5),
    @("trainingaudio/15_kopota3.wav", "pngimages/15_barrel.png", "trainingaudio/04_kitoti2.wav", "pngimages/04_ladder.png", -0.5, 0.5),
    @("trainingaudio/07_pitapi2.wav", "pngimages/07_suitcase.png", "trainingaudio/06_titoka3.wav", "pngimages/06_tent.png", 0.5, -0.5),
    @("trainingaudio/20_tatito1.wav", "pngimages/20_pizza.png", "trainingaudio/01_kitipi1.wav", "pngimages/01_gift.png", -0.5, 0.5),
    @("trainingaudio/07_pitapi2.wav", "pngimages/07_suitcase.png", "trainingaudio/21_papika1.wav", "pngimages/21_cheese.png", 0.5, -0.5),
    @("trainingaudio/04_kitoti2.wav", "pngimages/04_ladder.png", "trainingaudio/18_popata2.wav", "pngimages/18_donut.png", -0.5, 0.5),
    @("trainingaudio/15_kopota3.wav", "pngimages/15_barrel.png", "trainingaudio/19_papipi1.wav", "pngimages/19_burger.png", -0.5, 0.5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
